# The document has three inline pictures living in the header/footer
# stories whose DrawingML "name" label (wp:docPr/@name, mirrored onto
# pic:cNvPr/@name) needs to be swapped, per the target diff:
#   - footer (first page)  id=3  Pearson logo: image2.png -> image1.png
#   - footer (default)     id=2  Pearson logo: image2.png -> image1.png
#   - header (first page)  id=1  BTec logo:    image1.jpg -> image2.jpg
#
# Everything else in the document (text, structure, rels, media) stays
# untouched.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# -- Footers: the Pearson Edexcel logo picture appears in both the
#    default footer and the first-page footer; both get the same
#    rename (image2.png -> image1.png).
for ($f = 1; $f -le 3; $f++) {
    $ftr = $sec.Footers.Item($f)
    if ($ftr.Exists) {
        for ($i = 1; $i -le $ftr.Range.InlineShapes.Count; $i++) {
            $shp = $ftr.Range.InlineShapes.Item($i)
            if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                $shp.Name = "image1.png"
            }
        }
    }
}

# -- Header: the BTEC logo picture lives in the first-page header
#    (image1.jpg -> image2.jpg).
for ($h = 1; $h -le 3; $h++) {
    $hdr = $sec.Headers.Item($h)
    if ($hdr.Exists) {
        for ($i = 1; $i -le $hdr.Range.InlineShapes.Count; $i++) {
            $shp = $hdr.Range.InlineShapes.Item($i)
            if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                $shp.Name = "image2.jpg"
            }
        }
    }
}
